# Add new account rows into the "Export" sheet, keeping the existing
# descending-by-Saldo ordering. Rows are inserted bottom-most first so
# that the row numbers referenced below (taken from the ORIGINAL sheet
# layout) stay valid while we work.
#
# Original layout (1-based rows):
#   1  Conta / Nome / Saldo      (header)
#   2  001882235 LAGO            278323.91
#   3  004368468 AHMAD           88330.09   <- new row inserted before this
#   4  004211368 ILTON           27558.5    <- new rows inserted before this
#   5  004222784 RAFAEL          13000      <- new rows inserted before this
#   6  004243043 SUELI           12752.5
#   7  004313254 GUSTAVO         4292       <- new row inserted before this

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-AccountRow($rowNum, $conta, $nome, $saldo) {
    $ws.Rows.Item($rowNum).Insert()
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $conta
    $ws.Cells.Item($rowNum, 2).Value = $nome
    $ws.Cells.Item($rowNum, 3).Value = $saldo
}

# Insert bottom-to-top so earlier original row indices remain correct.

# Before original row 7 (GUSTAVO) -> i.e. right after SUELI (12752.5)
Add-AccountRow 7 "005995120" "ERIK" 10069.58

# Before original row 5 (RAFAEL)
Add-AccountRow 5 "005701765" "F" 17118.29
Add-AccountRow 6 "005018038" "ELAINE" 16537.23

# Before original row 4 (ILTON)
Add-AccountRow 4 "004386464" "CARLOS" 70000
Add-AccountRow 5 "004458624" "PEDRO" 49399.46
Add-AccountRow 6 "004451978" "ANTONIO" 33074.45

# Before original row 3 (AHMAD)
Add-AccountRow 3 "005637820" "GUILHERME" 210000
